# feat: add 2022-Q4 data
#
# Result:
#   - "总计"     : existing summary row becomes 2022-Q4, a new row below it
#                  keeps the old 2022-Q3 figures.
#   - "2022-Q4" : new sheet (takes over the slot/rId formerly used by
#                  "2022-Q3") with the latest fund-position numbers.
#   - "2022-Q3" : new sheet holding exactly what used to be there, moved
#                  to the end.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Row 2 now reports 2022-Q4 (count / market value unchanged)
$total.Range("B2").Value = "2022-Q4"

# New row 3 carries what used to be the 2022-Q3 figures
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.14

# Give A3 the same look as A2 (bold / centered / bordered header style)
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: the sheet currently named "2022-Q3" is first duplicated (so the
# duplicate preserves the old 2022-Q3 data verbatim), then the original
# is renamed/repurposed to hold the new 2022-Q4 numbers, and the
# duplicate is renamed back to "2022-Q3".
# ---------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)
$oldQ3.Copy($null, $oldQ3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$newQ3 = $wb.Worksheets.Item(3)
$newQ3.Name = "2022-Q3"

# $q4 is a brand-new sheet in spirit, so it gets the same default page
# margins new sheets in this workbook use (i.e. the same as "总计"),
# rather than the margins inherited from the old "2022-Q3" copy.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Step 3: overwrite $q4 with the latest fund-position figures
# ---------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0

# These look numeric but must be stored as text, matching the source data
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "004685"
$q4.Range("B2").ClearFormats()

$q4.Range("C2").Value = "金元顺安元启灵活配置混合"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "15.29"
$q4.Range("D2").ClearFormats()

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "76.11"
$q4.Range("E2").ClearFormats()

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "0.92"
$q4.Range("F2").ClearFormats()

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1407"
$q4.Range("G2").ClearFormats()

$q4.Range("H2").Value = 9

# Re-apply the header/index styling (ClearFormats above stripped it)
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)

# Keep "总计" as the active sheet (matches the original workbook state)
$total.Activate()
